$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: replace "Papeleros" event with "Musica clasica" event ---
$ws.Range("A2").Value = "Musica clasica"
$ws.Range("B2").Value = "Concierto musica clásica"
$ws.Range("C2").Value = 45902.5
$ws.Range("D2").Value = "Media torta "
$ws.Range("E2").Value = "Av. 26 15 - 15"
$ws.Range("F2").Value = "Bogotá"
$ws.Range("G2").Value = 4.711
$ws.Range("H2").Value = -74.0421

# Row 2 shrinks from height 45 to 30
$ws.Rows.Item(2).RowHeight = 30

# --- Row 3: fill in new "Musica Instrumental" event (previously blank) ---
# Bring G3/H3 formatting in line with G2/H2 (general number format instead
# of the default text format) before writing numeric values into them.
$ws.Range("G2").Copy()
$ws.Range("G3").PasteSpecial(-4122)
$ws.Range("H2").Copy()
$ws.Range("H3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A3").Value = "Musica Instrumental"
$ws.Range("B3").Value = "Concierto musica Peruana"
$ws.Range("C3").Value = 45916.5
$ws.Range("D3").Value = "Centro de convenciones Lima"
$ws.Range("E3").Value = "Dg. 140 15 - 40"
$ws.Range("F3").Value = "Lima"
$ws.Range("G3").Value = -12.1211
$ws.Range("H3").Value = -77.0297

# Row 3 takes the same height as row 2
$ws.Rows.Item(3).RowHeight = 30
